$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) to snake_case names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the Spanish connector words ("de", "del", "el", "la", "las", "los", "y") ---
# in specific municipality / state names throughout the data rows
$ws.Range("B9").Value   = "Mazapa De Madero"
$ws.Range("B11").Value  = "San Cristóbal De Las Casas"
$ws.Range("A22").Value  = "Ciudad De México"
$ws.Range("A29").Value  = "Estado De México"
$ws.Range("B40").Value  = "Jaral Del Progreso"
$ws.Range("B44").Value  = "Purísima Del Rincón"
$ws.Range("B50").Value  = "Coyuca De Catalán"
$ws.Range("B52").Value  = "Huitzuco De Los Figueroa"
$ws.Range("B55").Value  = "Técpan De Galeana"
$ws.Range("B56").Value  = "Tlapa De Comonfort"
$ws.Range("B64").Value  = "Tenango De Doria"
$ws.Range("B65").Value  = "Tulancingo De Bravo"
$ws.Range("B66").Value  = "Zacualtipán De Ángeles"
$ws.Range("B69").Value  = "Atotonilco El Alto"
$ws.Range("B70").Value  = "Autlán De Navarro"
$ws.Range("B72").Value  = "Encarnación De Díaz"
$ws.Range("B89").Value  = "Mier Y Noriega"
$ws.Range("B90").Value  = "San Nicolás De Los Garza"
$ws.Range("B92").Value  = "Acatlán De Pérez Figueroa"
$ws.Range("B93").Value  = "Oaxaca De Juárez"
$ws.Range("B94").Value  = "San Francisco Del Mar"
$ws.Range("B101").Value = "Santo Domingo De Morelos"
$ws.Range("B103").Value = "Teotitlán De Flores Magón"
$ws.Range("B104").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B105").Value = "Tlacolula De Matamoros"
$ws.Range("B110").Value = "Amealco De Bonfil"
$ws.Range("B111").Value = "Jalpan De Serra"
$ws.Range("B112").Value = "Pinal De Amoles"
$ws.Range("B117").Value = "San Ciro De Acosta"
$ws.Range("B126").Value = "Hueyapan De Ocampo"
$ws.Range("B128").Value = "Ixhuatlán Del Café"

# --- Remove trailing metadata / duplicated footer rows (144-148 and 476-480) ---
# and anything beyond row 142, collapsing the used range back down to A1:D142
$ws.Rows("144:480").Delete()
